# Append 45 new "master-reg_center_device_h" rows (rows 102-146), matching
# the existing table's pattern: regcntr_id cycles 10002..10010, device_id
# increments by 1 starting at 3000121, and the remaining columns repeat the
# same constant values used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$regCenterIds = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)

$startRow = 102
$endRow = 146
$deviceId = 3000121

for ($row = $startRow; $row -le $endRow; $row++) {
    $idx = $row - $startRow
    $regCntrId = $regCenterIds[$idx % $regCenterIds.Length]

    $ws.Cells.Item($row, 1).Value = $regCntrId
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"

    $deviceId++
}

# Match the updated view state: selection covering the newly added rows.
$ws.Range("A102:B146").Select()

# Match the page setup change (portrait orientation).
$ws.PageSetup.Orientation = 1
